$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 33 (shifts existing rows 33:55 down to 34:56)
$ws.Rows(33).Insert()

# Copy formatting (thin-border style) from the row below into the new row,
# so the new row reuses the existing style instead of creating a new one.
$ws.Range("A34:F34").Copy()
$ws.Range("A33:F33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate the new product row: 20141201 / WALLS BALL GRAPE,
# belonging to group "4" position "11" (right after 20141200 WALLS BALL STRW GUAV).
$ws.Range("A33").Value = "20141201"
$ws.Range("B33").Value = "WALLS BALL GRAPE"
$ws.Range("C33").Value = "RWS05"
$ws.Range("D33").Value = "4"
$ws.Range("E33").Value = "11"
$ws.Range("F33").Value = "RT,(E-1B)"
